# Atualização de bases das ligas, do dia: 29-03-2024 às 17:05
#
# This script updates fixture rows that were re-ordered / refreshed on the
# source site: rows 118-120 (Clausura fixtures) are cyclically rotated, and
# rows 161-167 (Apertura fixtures) are shifted up by one (the stale row with
# id 7994680 is dropped, a new row of odds appears at the end of the block,
# and the surviving rows pick up their refreshed closing odds) which results
# in the last row (167) being removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: columns B..AC hold the match data for each row (column A is the
# running id and is left untouched).
# ---------------------------------------------------------------------------
$dataCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Get-RowValues($row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value2
    }
    return $vals
}

function Set-RowValues($row, $cols, $vals) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = $vals[$c]
    }
}

# ---------------------------------------------------------------------------
# 1) Rows 118-120: cyclic rotation of the match data (B:AC). Row 118 takes
#    row 119's data, row 119 takes row 120's data, and row 120 takes row
#    118's (original) data. The running id in column A stays as-is.
# ---------------------------------------------------------------------------
$row118 = Get-RowValues 118 $dataCols
$row119 = Get-RowValues 119 $dataCols
$row120 = Get-RowValues 120 $dataCols

Set-RowValues 118 $dataCols $row119
Set-RowValues 119 $dataCols $row120
Set-RowValues 120 $dataCols $row118

# The freshly refreshed odds feed reported AB120 (PL_Ahh) as -1 instead of
# the value that would come purely from the rotation (-0.5).
$ws.Range("AB120").Value = -1

# ---------------------------------------------------------------------------
# 2) Rows 161-167 (fixtures without results yet, columns B:AA only): the
#    whole block shifts up by one row - row 161 becomes the old row 162,
#    ... , row 166 becomes the old row 167 - and the resulting duplicate
#    last row is removed.
# ---------------------------------------------------------------------------
$blockCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

$blockRows = @{}
for ($r = 161; $r -le 167; $r++) {
    $blockRows[$r] = Get-RowValues $r $blockCols
}

for ($r = 161; $r -le 166; $r++) {
    Set-RowValues $r $blockCols $blockRows[$r + 1]
}

# Updated closing odds (oddH/oddD/oddA) picked up by the odds feed for the
# two matches that moved into rows 161 and 164.
$ws.Range("U161").Value = 2.025
$ws.Range("V161").Value = 1.825

$ws.Range("R164").Value = 2.025
$ws.Range("S164").Value = 1.825
$ws.Range("U164").Value = 1.925
$ws.Range("V164").Value = 1.925

# Remove the now-duplicated trailing row, shrinking the used range to AC166.
$ws.Rows(167).Delete()
